$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Projeto 4"
$ws.Range("B5").Value = "Cliente 4"
$ws.Range("C5").Value = "Mensagem de exemplo 4"
$ws.Range("D5").Value = "Web, Mobile & Software"
$ws.Range("E5").Value = "url 4"
